# issue #5: stock data output to json file
# Adds a "property_category" column to the 股票 (stock) sheet, populated
# with the literal value "stock" for every data row, inserted right after
# the "total" column and before the "date" column (columns shift right).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (pushes old H:J date/legislator_name/legislator_id
# to I:K, copying formatting from the column that used to be there).
$ws.Columns.Item(8).Insert()

# Header for the freshly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Populate every stock data row with the new category value.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
